$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; old C "Image_id" and D "class_type" shift right
# (contents row2_b..row2_e.PNG move from C3:C6 to D3:D6 automatically)
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = 32.33

# New content, entered in the same order the author would have typed it so the
# shared-string table is rebuilt in the same sequence.
$ws.Range("E16").Value = "A"
$ws.Range("C1").Value = "Accident Id"
$ws.Range("C2").Value = "A-2827637"
$ws.Range("D1").Value = "Image_link   "
$ws.Range("D10").Value = "                  "
$ws.Range("D2").Value = "00001.PNG"

$ws.Range("D3").Select()
